# Auto-generated edit script: adds 2022-Q3 sheet + summary row
$wb = $excel.ActiveWorkbook

# ---- Style templates (reuse existing style indices instead of inventing new ones) ----
# style "2" (bold + border + centered) lives on existing header cells, e.g. summary sheet B1
$summary = $wb.Worksheets.Item(1)
$styleBold = $summary.Cells.Item(1,2)

# ================= 1. Insert the new '2022-Q3' worksheet right after the summary sheet =================
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = '2022-Q3'

# -- header row --
$styleBold.Copy($newSheet.Cells.Item(1,2))
$newSheet.Cells.Item(1,2).Value2 = '基金代码'
$styleBold.Copy($newSheet.Cells.Item(1,3))
$newSheet.Cells.Item(1,3).Value2 = '基金名称'
$styleBold.Copy($newSheet.Cells.Item(1,4))
$newSheet.Cells.Item(1,4).Value2 = '基金规模'
$styleBold.Copy($newSheet.Cells.Item(1,5))
$newSheet.Cells.Item(1,5).Value2 = '股票总仓位'
$styleBold.Copy($newSheet.Cells.Item(1,6))
$newSheet.Cells.Item(1,6).Value2 = '仓位占比'
$styleBold.Copy($newSheet.Cells.Item(1,7))
$newSheet.Cells.Item(1,7).Value2 = '持有市值(亿元)'
$styleBold.Copy($newSheet.Cells.Item(1,8))
$newSheet.Cells.Item(1,8).Value2 = '仓位排名'

# -- data rows (row 2 .. row 34) --
# row 2
$styleBold.Copy($newSheet.Cells.Item(2,1))
$newSheet.Cells.Item(2,1).Value2 = 0
$newSheet.Cells.Item(2,2).NumberFormat = "@"
$newSheet.Cells.Item(2,2).Value2 = '002685'
$newSheet.Cells.Item(2,3).Value2 = '中欧丰泓沪港深灵活配置混合A'
$newSheet.Cells.Item(2,4).NumberFormat = "@"
$newSheet.Cells.Item(2,4).Value2 = '42.36'
$newSheet.Cells.Item(2,5).NumberFormat = "@"
$newSheet.Cells.Item(2,5).Value2 = '92.77'
$newSheet.Cells.Item(2,6).NumberFormat = "@"
$newSheet.Cells.Item(2,6).Value2 = '7.67'
$newSheet.Cells.Item(2,7).NumberFormat = "@"
$newSheet.Cells.Item(2,7).Value2 = '3.2490'
$newSheet.Cells.Item(2,8).Value2 = 3
# row 3
$styleBold.Copy($newSheet.Cells.Item(3,1))
$newSheet.Cells.Item(3,1).Value2 = 1
$newSheet.Cells.Item(3,2).NumberFormat = "@"
$newSheet.Cells.Item(3,2).Value2 = '005847'
$newSheet.Cells.Item(3,3).Value2 = '富国沪港深业绩驱动混合A'
$newSheet.Cells.Item(3,4).NumberFormat = "@"
$newSheet.Cells.Item(3,4).Value2 = '32.84'
$newSheet.Cells.Item(3,5).NumberFormat = "@"
$newSheet.Cells.Item(3,5).Value2 = '85.47'
$newSheet.Cells.Item(3,6).NumberFormat = "@"
$newSheet.Cells.Item(3,6).Value2 = '8.46'
$newSheet.Cells.Item(3,7).NumberFormat = "@"
$newSheet.Cells.Item(3,7).Value2 = '2.7783'
$newSheet.Cells.Item(3,8).Value2 = 3
# row 4
$styleBold.Copy($newSheet.Cells.Item(4,1))
$newSheet.Cells.Item(4,1).Value2 = 2
$newSheet.Cells.Item(4,2).NumberFormat = "@"
$newSheet.Cells.Item(4,2).Value2 = '501087'
$newSheet.Cells.Item(4,3).Value2 = '交银施罗德瑞丰混合（LOF）'
$newSheet.Cells.Item(4,4).NumberFormat = "@"
$newSheet.Cells.Item(4,4).Value2 = '23.31'
$newSheet.Cells.Item(4,5).NumberFormat = "@"
$newSheet.Cells.Item(4,5).Value2 = '84.93'
$newSheet.Cells.Item(4,6).NumberFormat = "@"
$newSheet.Cells.Item(4,6).Value2 = '7.89'
$newSheet.Cells.Item(4,7).NumberFormat = "@"
$newSheet.Cells.Item(4,7).Value2 = '1.8392'
$newSheet.Cells.Item(4,8).Value2 = 5
# row 5
$styleBold.Copy($newSheet.Cells.Item(5,1))
$newSheet.Cells.Item(5,1).Value2 = 3
$newSheet.Cells.Item(5,2).NumberFormat = "@"
$newSheet.Cells.Item(5,2).Value2 = '007455'
$newSheet.Cells.Item(5,3).Value2 = '富国蓝筹精选股票（QDII）人民币'
$newSheet.Cells.Item(5,4).NumberFormat = "@"
$newSheet.Cells.Item(5,4).Value2 = '13.62'
$newSheet.Cells.Item(5,5).NumberFormat = "@"
$newSheet.Cells.Item(5,5).Value2 = '85.59'
$newSheet.Cells.Item(5,6).NumberFormat = "@"
$newSheet.Cells.Item(5,6).Value2 = '7.95'
$newSheet.Cells.Item(5,7).NumberFormat = "@"
$newSheet.Cells.Item(5,7).Value2 = '1.0828'
$newSheet.Cells.Item(5,8).Value2 = 1
# row 6
$styleBold.Copy($newSheet.Cells.Item(6,1))
$newSheet.Cells.Item(6,1).Value2 = 4
$newSheet.Cells.Item(6,2).NumberFormat = "@"
$newSheet.Cells.Item(6,2).Value2 = '010583'
$newSheet.Cells.Item(6,3).Value2 = '富国蓝筹精选股票（QDII）美元'
$newSheet.Cells.Item(6,4).NumberFormat = "@"
$newSheet.Cells.Item(6,4).Value2 = '13.62'
$newSheet.Cells.Item(6,5).NumberFormat = "@"
$newSheet.Cells.Item(6,5).Value2 = '85.59'
$newSheet.Cells.Item(6,6).NumberFormat = "@"
$newSheet.Cells.Item(6,6).Value2 = '7.95'
$newSheet.Cells.Item(6,7).NumberFormat = "@"
$newSheet.Cells.Item(6,7).Value2 = '1.0828'
$newSheet.Cells.Item(6,8).Value2 = 1
# row 7
$styleBold.Copy($newSheet.Cells.Item(7,1))
$newSheet.Cells.Item(7,1).Value2 = 5
$newSheet.Cells.Item(7,2).NumberFormat = "@"
$newSheet.Cells.Item(7,2).Value2 = '001605'
$newSheet.Cells.Item(7,3).Value2 = '国富沪港深成长精选股票'
$newSheet.Cells.Item(7,4).NumberFormat = "@"
$newSheet.Cells.Item(7,4).Value2 = '27.46'
$newSheet.Cells.Item(7,5).NumberFormat = "@"
$newSheet.Cells.Item(7,5).Value2 = '84.27'
$newSheet.Cells.Item(7,6).NumberFormat = "@"
$newSheet.Cells.Item(7,6).Value2 = '3.23'
$newSheet.Cells.Item(7,7).NumberFormat = "@"
$newSheet.Cells.Item(7,7).Value2 = '0.8870'
$newSheet.Cells.Item(7,8).Value2 = 3
# row 8
$styleBold.Copy($newSheet.Cells.Item(8,1))
$newSheet.Cells.Item(8,1).Value2 = 6
$newSheet.Cells.Item(8,2).NumberFormat = "@"
$newSheet.Cells.Item(8,2).Value2 = '000934'
$newSheet.Cells.Item(8,3).Value2 = '国富大中华精选混合（QDII）'
$newSheet.Cells.Item(8,4).NumberFormat = "@"
$newSheet.Cells.Item(8,4).Value2 = '19.83'
$newSheet.Cells.Item(8,5).NumberFormat = "@"
$newSheet.Cells.Item(8,5).Value2 = '72.45'
$newSheet.Cells.Item(8,6).NumberFormat = "@"
$newSheet.Cells.Item(8,6).Value2 = '3.26'
$newSheet.Cells.Item(8,7).NumberFormat = "@"
$newSheet.Cells.Item(8,7).Value2 = '0.6465'
$newSheet.Cells.Item(8,8).Value2 = 3
# row 9
$styleBold.Copy($newSheet.Cells.Item(9,1))
$newSheet.Cells.Item(9,1).Value2 = 7
$newSheet.Cells.Item(9,2).NumberFormat = "@"
$newSheet.Cells.Item(9,2).Value2 = '006370'
$newSheet.Cells.Item(9,3).Value2 = '国富大中华精选混合（QDII）美元'
$newSheet.Cells.Item(9,4).NumberFormat = "@"
$newSheet.Cells.Item(9,4).Value2 = '19.83'
$newSheet.Cells.Item(9,5).NumberFormat = "@"
$newSheet.Cells.Item(9,5).Value2 = '72.45'
$newSheet.Cells.Item(9,6).NumberFormat = "@"
$newSheet.Cells.Item(9,6).Value2 = '3.26'
$newSheet.Cells.Item(9,7).NumberFormat = "@"
$newSheet.Cells.Item(9,7).Value2 = '0.6465'
$newSheet.Cells.Item(9,8).Value2 = 3
# row 10
$styleBold.Copy($newSheet.Cells.Item(10,1))
$newSheet.Cells.Item(10,1).Value2 = 8
$newSheet.Cells.Item(10,2).NumberFormat = "@"
$newSheet.Cells.Item(10,2).Value2 = '002686'
$newSheet.Cells.Item(10,3).Value2 = '中欧丰泓沪港深灵活配置混合C'
$newSheet.Cells.Item(10,4).NumberFormat = "@"
$newSheet.Cells.Item(10,4).Value2 = '7.40'
$newSheet.Cells.Item(10,5).NumberFormat = "@"
$newSheet.Cells.Item(10,5).Value2 = '92.77'
$newSheet.Cells.Item(10,6).NumberFormat = "@"
$newSheet.Cells.Item(10,6).Value2 = '7.67'
$newSheet.Cells.Item(10,7).NumberFormat = "@"
$newSheet.Cells.Item(10,7).Value2 = '0.5676'
$newSheet.Cells.Item(10,8).Value2 = 3
# row 11
$styleBold.Copy($newSheet.Cells.Item(11,1))
$newSheet.Cells.Item(11,1).Value2 = 9
$newSheet.Cells.Item(11,2).NumberFormat = "@"
$newSheet.Cells.Item(11,2).Value2 = '009846'
$newSheet.Cells.Item(11,3).Value2 = '富兰克林国海港股通远见价值混合'
$newSheet.Cells.Item(11,4).NumberFormat = "@"
$newSheet.Cells.Item(11,4).Value2 = '12.78'
$newSheet.Cells.Item(11,5).NumberFormat = "@"
$newSheet.Cells.Item(11,5).Value2 = '83.81'
$newSheet.Cells.Item(11,6).NumberFormat = "@"
$newSheet.Cells.Item(11,6).Value2 = '3.70'
$newSheet.Cells.Item(11,7).NumberFormat = "@"
$newSheet.Cells.Item(11,7).Value2 = '0.4729'
$newSheet.Cells.Item(11,8).Value2 = 4
# row 12
$styleBold.Copy($newSheet.Cells.Item(12,1))
$newSheet.Cells.Item(12,1).Value2 = 10
$newSheet.Cells.Item(12,2).NumberFormat = "@"
$newSheet.Cells.Item(12,2).Value2 = '013991'
$newSheet.Cells.Item(12,3).Value2 = '中欧港股通精选一年持有混合A'
$newSheet.Cells.Item(12,4).NumberFormat = "@"
$newSheet.Cells.Item(12,4).Value2 = '6.69'
$newSheet.Cells.Item(12,5).NumberFormat = "@"
$newSheet.Cells.Item(12,5).Value2 = '93.38'
$newSheet.Cells.Item(12,6).NumberFormat = "@"
$newSheet.Cells.Item(12,6).Value2 = '6.46'
$newSheet.Cells.Item(12,7).NumberFormat = "@"
$newSheet.Cells.Item(12,7).Value2 = '0.4322'
$newSheet.Cells.Item(12,8).Value2 = 4
# row 13
$styleBold.Copy($newSheet.Cells.Item(13,1))
$newSheet.Cells.Item(13,1).Value2 = 11
$newSheet.Cells.Item(13,2).NumberFormat = "@"
$newSheet.Cells.Item(13,2).Value2 = '011635'
$newSheet.Cells.Item(13,3).Value2 = '富国港股通策略精选混合A'
$newSheet.Cells.Item(13,4).NumberFormat = "@"
$newSheet.Cells.Item(13,4).Value2 = '6.21'
$newSheet.Cells.Item(13,5).NumberFormat = "@"
$newSheet.Cells.Item(13,5).Value2 = '73.36'
$newSheet.Cells.Item(13,6).NumberFormat = "@"
$newSheet.Cells.Item(13,6).Value2 = '5.86'
$newSheet.Cells.Item(13,7).NumberFormat = "@"
$newSheet.Cells.Item(13,7).Value2 = '0.3639'
$newSheet.Cells.Item(13,8).Value2 = 1
# row 14
$styleBold.Copy($newSheet.Cells.Item(14,1))
$newSheet.Cells.Item(14,1).Value2 = 12
$newSheet.Cells.Item(14,2).NumberFormat = "@"
$newSheet.Cells.Item(14,2).Value2 = '012744'
$newSheet.Cells.Item(14,3).Value2 = '光大保德信品质生活混合A'
$newSheet.Cells.Item(14,4).NumberFormat = "@"
$newSheet.Cells.Item(14,4).Value2 = '5.60'
$newSheet.Cells.Item(14,5).NumberFormat = "@"
$newSheet.Cells.Item(14,5).Value2 = '88.62'
$newSheet.Cells.Item(14,6).NumberFormat = "@"
$newSheet.Cells.Item(14,6).Value2 = '5.46'
$newSheet.Cells.Item(14,7).NumberFormat = "@"
$newSheet.Cells.Item(14,7).Value2 = '0.3058'
$newSheet.Cells.Item(14,8).Value2 = 6
# row 15
$styleBold.Copy($newSheet.Cells.Item(15,1))
$newSheet.Cells.Item(15,1).Value2 = 13
$newSheet.Cells.Item(15,2).NumberFormat = "@"
$newSheet.Cells.Item(15,2).Value2 = '013992'
$newSheet.Cells.Item(15,3).Value2 = '中欧港股通精选一年持有混合C'
$newSheet.Cells.Item(15,4).NumberFormat = "@"
$newSheet.Cells.Item(15,4).Value2 = '4.68'
$newSheet.Cells.Item(15,5).NumberFormat = "@"
$newSheet.Cells.Item(15,5).Value2 = '93.38'
$newSheet.Cells.Item(15,6).NumberFormat = "@"
$newSheet.Cells.Item(15,6).Value2 = '6.46'
$newSheet.Cells.Item(15,7).NumberFormat = "@"
$newSheet.Cells.Item(15,7).Value2 = '0.3023'
$newSheet.Cells.Item(15,8).Value2 = 4
# row 16
$styleBold.Copy($newSheet.Cells.Item(16,1))
$newSheet.Cells.Item(16,1).Value2 = 14
$newSheet.Cells.Item(16,2).NumberFormat = "@"
$newSheet.Cells.Item(16,2).Value2 = '011117'
$newSheet.Cells.Item(16,3).Value2 = '富国沪港深业绩驱动混合C'
$newSheet.Cells.Item(16,4).NumberFormat = "@"
$newSheet.Cells.Item(16,4).Value2 = '3.53'
$newSheet.Cells.Item(16,5).NumberFormat = "@"
$newSheet.Cells.Item(16,5).Value2 = '85.47'
$newSheet.Cells.Item(16,6).NumberFormat = "@"
$newSheet.Cells.Item(16,6).Value2 = '8.46'
$newSheet.Cells.Item(16,7).NumberFormat = "@"
$newSheet.Cells.Item(16,7).Value2 = '0.2986'
$newSheet.Cells.Item(16,8).Value2 = 3
# row 17
$styleBold.Copy($newSheet.Cells.Item(17,1))
$newSheet.Cells.Item(17,1).Value2 = 15
$newSheet.Cells.Item(17,2).NumberFormat = "@"
$newSheet.Cells.Item(17,2).Value2 = '010088'
$newSheet.Cells.Item(17,3).Value2 = '工银优质成长混合A'
$newSheet.Cells.Item(17,4).NumberFormat = "@"
$newSheet.Cells.Item(17,4).Value2 = '15.38'
$newSheet.Cells.Item(17,5).NumberFormat = "@"
$newSheet.Cells.Item(17,5).Value2 = '69.60'
$newSheet.Cells.Item(17,6).NumberFormat = "@"
$newSheet.Cells.Item(17,6).Value2 = '1.88'
$newSheet.Cells.Item(17,7).NumberFormat = "@"
$newSheet.Cells.Item(17,7).Value2 = '0.2891'
$newSheet.Cells.Item(17,8).Value2 = 9
# row 18
$styleBold.Copy($newSheet.Cells.Item(18,1))
$newSheet.Cells.Item(18,1).Value2 = 16
$newSheet.Cells.Item(18,2).NumberFormat = "@"
$newSheet.Cells.Item(18,2).Value2 = '006039'
$newSheet.Cells.Item(18,3).Value2 = '国富估值优势混合'
$newSheet.Cells.Item(18,4).NumberFormat = "@"
$newSheet.Cells.Item(18,4).Value2 = '6.13'
$newSheet.Cells.Item(18,5).NumberFormat = "@"
$newSheet.Cells.Item(18,5).Value2 = '81.55'
$newSheet.Cells.Item(18,6).NumberFormat = "@"
$newSheet.Cells.Item(18,6).Value2 = '4.65'
$newSheet.Cells.Item(18,7).NumberFormat = "@"
$newSheet.Cells.Item(18,7).Value2 = '0.2850'
$newSheet.Cells.Item(18,8).Value2 = 1
# row 19
$styleBold.Copy($newSheet.Cells.Item(19,1))
$newSheet.Cells.Item(19,1).Value2 = 17
$newSheet.Cells.Item(19,2).NumberFormat = "@"
$newSheet.Cells.Item(19,2).Value2 = '012584'
$newSheet.Cells.Item(19,3).Value2 = '南方中国新兴经济9个月持有期混合（QDII）A'
$newSheet.Cells.Item(19,4).NumberFormat = "@"
$newSheet.Cells.Item(19,4).Value2 = '2.69'
$newSheet.Cells.Item(19,5).NumberFormat = "@"
$newSheet.Cells.Item(19,5).Value2 = '91.51'
$newSheet.Cells.Item(19,6).NumberFormat = "@"
$newSheet.Cells.Item(19,6).Value2 = '4.62'
$newSheet.Cells.Item(19,7).NumberFormat = "@"
$newSheet.Cells.Item(19,7).Value2 = '0.1243'
$newSheet.Cells.Item(19,8).Value2 = 4
# row 20
$styleBold.Copy($newSheet.Cells.Item(20,1))
$newSheet.Cells.Item(20,1).Value2 = 18
$newSheet.Cells.Item(20,2).NumberFormat = "@"
$newSheet.Cells.Item(20,2).Value2 = '457001'
$newSheet.Cells.Item(20,3).Value2 = '国富亚洲机会股票（QDII）'
$newSheet.Cells.Item(20,4).NumberFormat = "@"
$newSheet.Cells.Item(20,4).Value2 = '3.80'
$newSheet.Cells.Item(20,5).NumberFormat = "@"
$newSheet.Cells.Item(20,5).Value2 = '83.80'
$newSheet.Cells.Item(20,6).NumberFormat = "@"
$newSheet.Cells.Item(20,6).Value2 = '3.19'
$newSheet.Cells.Item(20,7).NumberFormat = "@"
$newSheet.Cells.Item(20,7).Value2 = '0.1212'
$newSheet.Cells.Item(20,8).Value2 = 4
# row 21
$styleBold.Copy($newSheet.Cells.Item(21,1))
$newSheet.Cells.Item(21,1).Value2 = 19
$newSheet.Cells.Item(21,2).NumberFormat = "@"
$newSheet.Cells.Item(21,2).Value2 = '160125'
$newSheet.Cells.Item(21,3).Value2 = '南方香港优选股票（QDII-LOF）'
$newSheet.Cells.Item(21,4).NumberFormat = "@"
$newSheet.Cells.Item(21,4).Value2 = '2.01'
$newSheet.Cells.Item(21,5).NumberFormat = "@"
$newSheet.Cells.Item(21,5).Value2 = '81.74'
$newSheet.Cells.Item(21,6).NumberFormat = "@"
$newSheet.Cells.Item(21,6).Value2 = '5.88'
$newSheet.Cells.Item(21,7).NumberFormat = "@"
$newSheet.Cells.Item(21,7).Value2 = '0.1182'
$newSheet.Cells.Item(21,8).Value2 = 1
# row 22
$styleBold.Copy($newSheet.Cells.Item(22,1))
$newSheet.Cells.Item(22,1).Value2 = 20
$newSheet.Cells.Item(22,2).NumberFormat = "@"
$newSheet.Cells.Item(22,2).Value2 = '012770'
$newSheet.Cells.Item(22,3).Value2 = '光大保德信创新生活混合'
$newSheet.Cells.Item(22,4).NumberFormat = "@"
$newSheet.Cells.Item(22,4).Value2 = '2.71'
$newSheet.Cells.Item(22,5).NumberFormat = "@"
$newSheet.Cells.Item(22,5).Value2 = '86.69'
$newSheet.Cells.Item(22,6).NumberFormat = "@"
$newSheet.Cells.Item(22,6).Value2 = '3.64'
$newSheet.Cells.Item(22,7).NumberFormat = "@"
$newSheet.Cells.Item(22,7).Value2 = '0.0986'
$newSheet.Cells.Item(22,8).Value2 = 7
# row 23
$styleBold.Copy($newSheet.Cells.Item(23,1))
$newSheet.Cells.Item(23,1).Value2 = 21
$newSheet.Cells.Item(23,2).NumberFormat = "@"
$newSheet.Cells.Item(23,2).Value2 = '011349'
$newSheet.Cells.Item(23,3).Value2 = '淳厚现代服务业股票A'
$newSheet.Cells.Item(23,4).NumberFormat = "@"
$newSheet.Cells.Item(23,4).Value2 = '2.55'
$newSheet.Cells.Item(23,5).NumberFormat = "@"
$newSheet.Cells.Item(23,5).Value2 = '79.82'
$newSheet.Cells.Item(23,6).NumberFormat = "@"
$newSheet.Cells.Item(23,6).Value2 = '2.89'
$newSheet.Cells.Item(23,7).NumberFormat = "@"
$newSheet.Cells.Item(23,7).Value2 = '0.0737'
$newSheet.Cells.Item(23,8).Value2 = 6
# row 24
$styleBold.Copy($newSheet.Cells.Item(24,1))
$newSheet.Cells.Item(24,1).Value2 = 22
$newSheet.Cells.Item(24,2).NumberFormat = "@"
$newSheet.Cells.Item(24,2).Value2 = '012060'
$newSheet.Cells.Item(24,3).Value2 = '富国全球消费精选混合（QDII）A'
$newSheet.Cells.Item(24,4).NumberFormat = "@"
$newSheet.Cells.Item(24,4).Value2 = '2.68'
$newSheet.Cells.Item(24,5).NumberFormat = "@"
$newSheet.Cells.Item(24,5).Value2 = '66.08'
$newSheet.Cells.Item(24,6).NumberFormat = "@"
$newSheet.Cells.Item(24,6).Value2 = '2.56'
$newSheet.Cells.Item(24,7).NumberFormat = "@"
$newSheet.Cells.Item(24,7).Value2 = '0.0686'
$newSheet.Cells.Item(24,8).Value2 = 7
# row 25
$styleBold.Copy($newSheet.Cells.Item(25,1))
$newSheet.Cells.Item(25,1).Value2 = 23
$newSheet.Cells.Item(25,2).NumberFormat = "@"
$newSheet.Cells.Item(25,2).Value2 = '012061'
$newSheet.Cells.Item(25,3).Value2 = '富国全球消费精选混合（QDII）美元现汇'
$newSheet.Cells.Item(25,4).NumberFormat = "@"
$newSheet.Cells.Item(25,4).Value2 = '2.68'
$newSheet.Cells.Item(25,5).NumberFormat = "@"
$newSheet.Cells.Item(25,5).Value2 = '66.08'
$newSheet.Cells.Item(25,6).NumberFormat = "@"
$newSheet.Cells.Item(25,6).Value2 = '2.56'
$newSheet.Cells.Item(25,7).NumberFormat = "@"
$newSheet.Cells.Item(25,7).Value2 = '0.0686'
$newSheet.Cells.Item(25,8).Value2 = 7
# row 26
$styleBold.Copy($newSheet.Cells.Item(26,1))
$newSheet.Cells.Item(26,1).Value2 = 24
$newSheet.Cells.Item(26,2).NumberFormat = "@"
$newSheet.Cells.Item(26,2).Value2 = '011636'
$newSheet.Cells.Item(26,3).Value2 = '富国港股通策略精选混合C'
$newSheet.Cells.Item(26,4).NumberFormat = "@"
$newSheet.Cells.Item(26,4).Value2 = '0.57'
$newSheet.Cells.Item(26,5).NumberFormat = "@"
$newSheet.Cells.Item(26,5).Value2 = '73.36'
$newSheet.Cells.Item(26,6).NumberFormat = "@"
$newSheet.Cells.Item(26,6).Value2 = '5.86'
$newSheet.Cells.Item(26,7).NumberFormat = "@"
$newSheet.Cells.Item(26,7).Value2 = '0.0334'
$newSheet.Cells.Item(26,8).Value2 = 1
# row 27
$styleBold.Copy($newSheet.Cells.Item(27,1))
$newSheet.Cells.Item(27,1).Value2 = 25
$newSheet.Cells.Item(27,2).NumberFormat = "@"
$newSheet.Cells.Item(27,2).Value2 = '010089'
$newSheet.Cells.Item(27,3).Value2 = '工银优质成长混合C'
$newSheet.Cells.Item(27,4).NumberFormat = "@"
$newSheet.Cells.Item(27,4).Value2 = '1.05'
$newSheet.Cells.Item(27,5).NumberFormat = "@"
$newSheet.Cells.Item(27,5).Value2 = '69.60'
$newSheet.Cells.Item(27,6).NumberFormat = "@"
$newSheet.Cells.Item(27,6).Value2 = '1.88'
$newSheet.Cells.Item(27,7).NumberFormat = "@"
$newSheet.Cells.Item(27,7).Value2 = '0.0197'
$newSheet.Cells.Item(27,8).Value2 = 9
# row 28
$styleBold.Copy($newSheet.Cells.Item(28,1))
$newSheet.Cells.Item(28,1).Value2 = 26
$newSheet.Cells.Item(28,2).NumberFormat = "@"
$newSheet.Cells.Item(28,2).Value2 = '012758'
$newSheet.Cells.Item(28,3).Value2 = '光大保德信品质生活混合C'
$newSheet.Cells.Item(28,4).NumberFormat = "@"
$newSheet.Cells.Item(28,4).Value2 = '0.35'
$newSheet.Cells.Item(28,5).NumberFormat = "@"
$newSheet.Cells.Item(28,5).Value2 = '88.62'
$newSheet.Cells.Item(28,6).NumberFormat = "@"
$newSheet.Cells.Item(28,6).Value2 = '5.46'
$newSheet.Cells.Item(28,7).NumberFormat = "@"
$newSheet.Cells.Item(28,7).Value2 = '0.0191'
$newSheet.Cells.Item(28,8).Value2 = 6
# row 29
$styleBold.Copy($newSheet.Cells.Item(29,1))
$newSheet.Cells.Item(29,1).Value2 = 27
$newSheet.Cells.Item(29,2).NumberFormat = "@"
$newSheet.Cells.Item(29,2).Value2 = '011350'
$newSheet.Cells.Item(29,3).Value2 = '淳厚现代服务业股票C'
$newSheet.Cells.Item(29,4).NumberFormat = "@"
$newSheet.Cells.Item(29,4).Value2 = '0.57'
$newSheet.Cells.Item(29,5).NumberFormat = "@"
$newSheet.Cells.Item(29,5).Value2 = '79.82'
$newSheet.Cells.Item(29,6).NumberFormat = "@"
$newSheet.Cells.Item(29,6).Value2 = '2.89'
$newSheet.Cells.Item(29,7).NumberFormat = "@"
$newSheet.Cells.Item(29,7).Value2 = '0.0165'
$newSheet.Cells.Item(29,8).Value2 = 6
# row 30
$styleBold.Copy($newSheet.Cells.Item(30,1))
$newSheet.Cells.Item(30,1).Value2 = 28
$newSheet.Cells.Item(30,2).NumberFormat = "@"
$newSheet.Cells.Item(30,2).Value2 = '006205'
$newSheet.Cells.Item(30,3).Value2 = '汇添富沪港深优势精选定期开放混合'
$newSheet.Cells.Item(30,4).NumberFormat = "@"
$newSheet.Cells.Item(30,4).Value2 = '0.33'
$newSheet.Cells.Item(30,5).NumberFormat = "@"
$newSheet.Cells.Item(30,5).Value2 = '95.18'
$newSheet.Cells.Item(30,6).NumberFormat = "@"
$newSheet.Cells.Item(30,6).Value2 = '4.28'
$newSheet.Cells.Item(30,7).NumberFormat = "@"
$newSheet.Cells.Item(30,7).Value2 = '0.0141'
$newSheet.Cells.Item(30,8).Value2 = 9
# row 31
$styleBold.Copy($newSheet.Cells.Item(31,1))
$newSheet.Cells.Item(31,1).Value2 = 29
$newSheet.Cells.Item(31,2).NumberFormat = "@"
$newSheet.Cells.Item(31,2).Value2 = '001942'
$newSheet.Cells.Item(31,3).Value2 = '前海开源沪港深汇鑫灵活配置混合A'
$newSheet.Cells.Item(31,4).NumberFormat = "@"
$newSheet.Cells.Item(31,4).Value2 = '0.17'
$newSheet.Cells.Item(31,5).NumberFormat = "@"
$newSheet.Cells.Item(31,5).Value2 = '87.24'
$newSheet.Cells.Item(31,6).NumberFormat = "@"
$newSheet.Cells.Item(31,6).Value2 = '4.59'
$newSheet.Cells.Item(31,7).NumberFormat = "@"
$newSheet.Cells.Item(31,7).Value2 = '0.0078'
$newSheet.Cells.Item(31,8).Value2 = 8
# row 32
$styleBold.Copy($newSheet.Cells.Item(32,1))
$newSheet.Cells.Item(32,1).Value2 = 30
$newSheet.Cells.Item(32,2).NumberFormat = "@"
$newSheet.Cells.Item(32,2).Value2 = '012585'
$newSheet.Cells.Item(32,3).Value2 = '南方中国新兴经济9个月持有期混合（QDII）C'
$newSheet.Cells.Item(32,4).NumberFormat = "@"
$newSheet.Cells.Item(32,4).Value2 = '0.10'
$newSheet.Cells.Item(32,5).NumberFormat = "@"
$newSheet.Cells.Item(32,5).Value2 = '91.51'
$newSheet.Cells.Item(32,6).NumberFormat = "@"
$newSheet.Cells.Item(32,6).Value2 = '4.62'
$newSheet.Cells.Item(32,7).NumberFormat = "@"
$newSheet.Cells.Item(32,7).Value2 = '0.0046'
$newSheet.Cells.Item(32,8).Value2 = 4
# row 33
$styleBold.Copy($newSheet.Cells.Item(33,1))
$newSheet.Cells.Item(33,1).Value2 = 31
$newSheet.Cells.Item(33,2).NumberFormat = "@"
$newSheet.Cells.Item(33,2).Value2 = '001943'
$newSheet.Cells.Item(33,3).Value2 = '前海开源沪港深汇鑫灵活配置混合C'
$newSheet.Cells.Item(33,4).NumberFormat = "@"
$newSheet.Cells.Item(33,4).Value2 = '0.09'
$newSheet.Cells.Item(33,5).NumberFormat = "@"
$newSheet.Cells.Item(33,5).Value2 = '87.24'
$newSheet.Cells.Item(33,6).NumberFormat = "@"
$newSheet.Cells.Item(33,6).Value2 = '4.59'
$newSheet.Cells.Item(33,7).NumberFormat = "@"
$newSheet.Cells.Item(33,7).Value2 = '0.0041'
$newSheet.Cells.Item(33,8).Value2 = 8
# row 34
$styleBold.Copy($newSheet.Cells.Item(34,1))
$newSheet.Cells.Item(34,1).Value2 = 32
$newSheet.Cells.Item(34,2).NumberFormat = "@"
$newSheet.Cells.Item(34,2).Value2 = '012062'
$newSheet.Cells.Item(34,3).Value2 = '富国全球消费精选混合（QDII）C'
$newSheet.Cells.Item(34,4).NumberFormat = "@"
$newSheet.Cells.Item(34,4).Value2 = '-2.54'
$newSheet.Cells.Item(34,5).NumberFormat = "@"
$newSheet.Cells.Item(34,5).Value2 = '66.08'
$newSheet.Cells.Item(34,6).NumberFormat = "@"
$newSheet.Cells.Item(34,6).Value2 = '2.56'
$newSheet.Cells.Item(34,7).NumberFormat = "@"
$newSheet.Cells.Item(34,7).Value2 = '-0.0650'
$newSheet.Cells.Item(34,8).Value2 = 7

# ================= 2. Update the summary ('总计') sheet: insert the 2022-Q3 row, shift the rest down =================
$styleBold.Copy($summary.Cells.Item(2,1))
$summary.Cells.Item(2,1).Value2 = 0
$summary.Cells.Item(2,2).Value2 = '2022-Q3'
$summary.Cells.Item(2,3).Value2 = 33
$summary.Cells.Item(2,4).Value2 = 16.26
$styleBold.Copy($summary.Cells.Item(3,1))
$summary.Cells.Item(3,1).Value2 = 1
$summary.Cells.Item(3,2).Value2 = '2022-Q2'
$summary.Cells.Item(3,3).Value2 = 19
$summary.Cells.Item(3,4).Value2 = 18.62
$styleBold.Copy($summary.Cells.Item(4,1))
$summary.Cells.Item(4,1).Value2 = 2
$summary.Cells.Item(4,2).Value2 = '2022-Q1'
$summary.Cells.Item(4,3).Value2 = 22
$summary.Cells.Item(4,4).Value2 = 20.3
$styleBold.Copy($summary.Cells.Item(5,1))
$summary.Cells.Item(5,1).Value2 = 3
$summary.Cells.Item(5,2).Value2 = '2021-Q4'
$summary.Cells.Item(5,3).Value2 = 12
$summary.Cells.Item(5,4).Value2 = 13.44
$styleBold.Copy($summary.Cells.Item(6,1))
$summary.Cells.Item(6,1).Value2 = 4
$summary.Cells.Item(6,2).Value2 = '2021-Q3'
$summary.Cells.Item(6,3).Value2 = 10
$summary.Cells.Item(6,4).Value2 = 13.53
$styleBold.Copy($summary.Cells.Item(7,1))
$summary.Cells.Item(7,1).Value2 = 5
$summary.Cells.Item(7,2).Value2 = '2021-Q2'
$summary.Cells.Item(7,3).Value2 = 20
$summary.Cells.Item(7,4).Value2 = 14.32
$styleBold.Copy($summary.Cells.Item(8,1))
$summary.Cells.Item(8,1).Value2 = 6
$summary.Cells.Item(8,2).Value2 = '2021-Q1'
$summary.Cells.Item(8,3).Value2 = 16
$summary.Cells.Item(8,4).Value2 = 6.89
